# Cancer type facet INS - update queries to use LIKE wildcard matching,
# case-insensitive ORDER BY, an extra Relative Citation Ratio case, and a
# revised "Data Location Details" expression for the Programs query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$programsQuery = "SELECT DISTINCT `r`n    prg.program_name AS `"Program`",`r`n    prg.website AS `"Website`",`r`n    prg.focus_area AS `"Focus Area`",`r`n    prg.cancer_type AS `"Cancer Type`",`r`n   CASE `r`n        WHEN prg.data_link IS NOT NULL THEN prg.website       `r`n        ELSE prg.data_link`r`n    END AS `"Data Location Details`"`r`nFROM `r`n    df_program prg`r`nWHERE `r`n     prg.cancer_type LIKE '%Esophageal Cancer%'`r`nORDER BY `r`n    lower(prg.program_name) ASC`r`nLIMIT 100;"

$statsQuery = "SELECT DISTINCT`r`n    COUNT(DISTINCT prg.program_id) AS `"Programs`",`r`n    COUNT(DISTINCT prj.project_id) AS `"Projects`",`r`n    COUNT(DISTINCT gnt.grant_id) AS `"Grants`",`r`n    COUNT(DISTINCT pub.pmid) AS `"Publications`"`r`nFROM `r`n    df_program prg`r`nLEFT JOIN `r`n    df_project prj ON prg.program_id = prj.`"program.program_id`"`r`nLEFT JOIN `r`n    df_grant gnt ON prj.project_id = gnt.`"project.project_id`"`r`nLEFT JOIN `r`n    df_publication pub ON prj.project_id = pub.`"project.project_id`"`r`nWHERE `r`n    prg.cancer_type LIKE '%Esophageal Cancer%';"

$projectsQuery = "SELECT DISTINCT`r`n    prj.project_id AS `"Project ID`", `r`n    prj.project_title AS `"Project Title`",`r`n    prj.org_name AS `"Organization`",`r`n    prj.project_start_date AS `"Project Start Date`",`r`n    prj.project_end_date AS `"Project End Date`"`r`nFROM `r`n    df_project prj`r`nLEFT JOIN `r`n    df_program prg ON prj.`"program.program_id`" = prg.program_id`r`nLEFT JOIN `r`n    df_grant gnt ON prj.project_id = gnt.`"project.project_id`"`r`nLEFT JOIN `r`n    df_publication pub ON prj.project_id = pub.`"project.project_id`"`r`nWHERE `r`n     prg.cancer_type LIKE '%Esophageal Cancer%'`r`nORDER BY `r`n    lower(prj.project_id) ASC`r`nLIMIT 100;"

$grantsQuery = "SELECT DISTINCT`r`n    gnt.grant_id AS `"Grant ID`", `r`n    prj.project_id AS `"Project`",`r`n    gnt.grant_title AS `"Grant Title`",`r`n    gnt.principal_investigators AS `"Principal Investigators`",`r`n    gnt.program_officers AS `"Program Officers`",`r`n    gnt.fiscal_year AS `"Fiscal Year`",`r`n    gnt.project_end_date AS `"Project End Date`"`r`nFROM `r`n    df_grant gnt`r`nLEFT JOIN `r`n    df_project prj ON gnt.`"project.project_id`" = prj.project_id`r`nLEFT JOIN `r`n    df_program prg ON prj.`"program.program_id`" = prg.program_id`r`nLEFT JOIN `r`n    df_publication pub ON prj.project_id = pub.`"project.project_id`"`r`nWHERE `r`n    prg.cancer_type LIKE '%Esophageal Cancer%'`r`nORDER BY `r`n    lower(gnt.grant_id) ASC`r`nLIMIT 100;"

$publicationsQuery = "SELECT DISTINCT`r`n    pub.pmid AS `"PubMed ID`", `r`n    pub.title AS `"Title`",`r`n    pub.authors AS `"Authors`",`r`n    pub.publication_date AS `"Publication Date`",`r`n    pub.cited_by AS `"Cited By`",`r`n    CASE `r`n    WHEN pub.relative_citation_ratio = 0 THEN '0'`r`n    WHEN pub.relative_citation_ratio = 7.0 THEN '7'`r`n    WHEN pub.relative_citation_ratio = 2.0 THEN '2'`r`n    WHEN pub.relative_citation_ratio = 3.0 THEN '3'`r`n    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) `r`n    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)`r`nEND AS `"Relative Citation Ratio`"`r`nFROM `r`n    df_publication pub`r`nLEFT JOIN `r`n    df_project prj ON pub.`"project.project_id`" = prj.project_id`r`nLEFT JOIN `r`n    df_program prg ON prj.`"program.program_id`" = prg.program_id`r`nLEFT JOIN `r`n    df_grant gnt ON prj.project_id = gnt.`"project.project_id`"`r`nWHERE `r`n     prg.cancer_type LIKE '%Esophageal Cancer%'`r`nORDER BY `r`n    lower(pub.pmid) ASC`r`nLIMIT 100;"

$ws.Range("B2").Value = $programsQuery
$ws.Range("C2").Value = $statsQuery
$ws.Range("B3").Value = $projectsQuery
$ws.Range("B4").Value = $grantsQuery
$ws.Range("B5").Value = $publicationsQuery

# Update the active selection on the sheet (was H2 with topLeftCell D1;
# now B2 with no explicit topLeftCell / frozen scroll position).
$null = $ws.Range("B2").Select()

# Best-effort: update workbook window position metadata.
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 12840
    $win.Top = 2440
} catch {
}
